# Updated symbol list on Sun Feb 12 06:32:54 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns on Sheet1 with the
# latest scrape values. Every cell in these two columns holds text
# (e.g. "307.95", "0.02%"), so we force a Text number format before
# writing the value -- otherwise Excel's COM layer would helpfully
# coerce numeric-looking strings (and percentages) into real numbers --
# then clear the temporary format again so the cell's formatting is
# left exactly as it was found (no stray per-cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "307.95"
Set-TextValue $ws.Range("E2") "0.02%"
Set-TextValue $ws.Range("D3") "41.02"
Set-TextValue $ws.Range("E3") "-0.24%"
Set-TextValue $ws.Range("D4") "5.237"
Set-TextValue $ws.Range("E4") "2.22%"
Set-TextValue $ws.Range("D5") "0.07673"
Set-TextValue $ws.Range("E5") "0.78%"
Set-TextValue $ws.Range("D6") "1.638"
Set-TextValue $ws.Range("E6") "1.17%"
Set-TextValue $ws.Range("D7") "0.9163"
Set-TextValue $ws.Range("E7") "1.92%"
Set-TextValue $ws.Range("D8") "2.441"
Set-TextValue $ws.Range("E8") "-0.29%"
Set-TextValue $ws.Range("D9") "0.1233"
Set-TextValue $ws.Range("E9") "13.33%"
Set-TextValue $ws.Range("D10") "0.1829"
Set-TextValue $ws.Range("E10") "3.54%"
Set-TextValue $ws.Range("D11") "0.09131"
Set-TextValue $ws.Range("E11") "-0.30%"
Set-TextValue $ws.Range("D12") "0.04264"
Set-TextValue $ws.Range("E12") "1.42%"
Set-TextValue $ws.Range("D13") "0.1051"
Set-TextValue $ws.Range("E13") "-0.03%"
Set-TextValue $ws.Range("D14") "0.001258"
Set-TextValue $ws.Range("E14") "0.70%"
Set-TextValue $ws.Range("D15") "0.005741"
Set-TextValue $ws.Range("E15") "-1.76%"
Set-TextValue $ws.Range("D17") "3.348"
Set-TextValue $ws.Range("E17") "-0.12%"
Set-TextValue $ws.Range("D18") "4.325"
Set-TextValue $ws.Range("E18") "1.67%"
Set-TextValue $ws.Range("D20") "7.306"
Set-TextValue $ws.Range("E20") "11.35%"
Set-TextValue $ws.Range("D21") "0.1385"
Set-TextValue $ws.Range("E21") "1.51%"
Set-TextValue $ws.Range("D22") "0.2895"
Set-TextValue $ws.Range("E22") "7.93%"
Set-TextValue $ws.Range("D23") "0.04081"
Set-TextValue $ws.Range("E23") "0.13%"
Set-TextValue $ws.Range("D24") "0.001264"
Set-TextValue $ws.Range("E24") "3.24%"
Set-TextValue $ws.Range("D25") "0.004328"
Set-TextValue $ws.Range("E25") "5.69%"
Set-TextValue $ws.Range("D26") "0.0001273"
Set-TextValue $ws.Range("E26") "-2.14%"
Set-TextValue $ws.Range("D38") "0.02466"
Set-TextValue $ws.Range("E38") "3.76%"
Set-TextValue $ws.Range("D39") "0.05298"
Set-TextValue $ws.Range("E39") "2.31%"
Set-TextValue $ws.Range("D40") "0.007854"
Set-TextValue $ws.Range("E40") "1.05%"
Set-TextValue $ws.Range("E41") "1.30%"
Set-TextValue $ws.Range("D42") "0.006793"
Set-TextValue $ws.Range("E42") "0.37%"
Set-TextValue $ws.Range("D43") "0.001915"
Set-TextValue $ws.Range("E43") "-1.88%"
Set-TextValue $ws.Range("D44") "0.007658"
Set-TextValue $ws.Range("E44") "-10.39%"
Set-TextValue $ws.Range("D45") "0.3062"
Set-TextValue $ws.Range("E45") "-0.28%"
Set-TextValue $ws.Range("E46") "-4.28%"
Set-TextValue $ws.Range("E47") "0.22%"
Set-TextValue $ws.Range("E48") "1,185.83%"
Set-TextValue $ws.Range("E49") "-26.03%"
Set-TextValue $ws.Range("E50") "0.22%"
Set-TextValue $ws.Range("E51") "0.22%"
